# south, yuktahar mistakes rectified
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Breakfast / Item 1 / Tuesday: "Poori" -> "Poha"
$ws.Range("D2").Value = "Poha"

# Lunch / Fry / Saturday: "Kakarakaya/Brinjal onions" -> "Kakarakaya"
$ws.Range("H7").Value = "Kakarakaya"

# Lunch / Spl rice / Saturday was blank -> add "Tomato Rice"
# (match the wrapped-text style already used by the rest of that column, e.g. I7)
$ws.Range("I7").Copy()
$ws.Range("H13").PasteSpecial(-4122)
$ws.Range("H13").Value = "Tomato Rice"

# Dinner / Fry / Monday: remove the erroneous "Aloo 65" entry
$ws.Range("C20").Value = ""

# leave the selection where the edits were last made
$ws.Range("C23").Select()
